# Auto-generated Excel COM-interop script
# Applies the numeric corrections described in the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2458.2
$ws.Range("I12").Value = 1763
$ws.Range("J12").Value = 3501
$ws.Range("K12").Value = 1763
$ws.Range("L12").Value = 3501
$ws.Range("M12").Value = -1593
$ws.Range("N12").Value = -3841

$ws.Range("H17").Value = 416375.88
$ws.Range("J17").Value = 473682.9
$ws.Range("L17").Value = 1421048.7
$ws.Range("N17").Value = -1421384.7

$ws.Range("H18").Value = 4824.5
$ws.Range("I18").Value = 6207.8335
$ws.Range("K18").Value = 6207.8335
$ws.Range("M18").Value = -5923.8335

$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H60").Value = 5000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H70").Value = 58336720
$ws.Range("J70").Value = 41669900
$ws.Range("L70").Value = 125009700
$ws.Range("N70").Value = -125010240

$ws.Range("H73").Value = 58336720
$ws.Range("J73").Value = 41669900
$ws.Range("L73").Value = 125009700
$ws.Range("N73").Value = -125011572

$ws.Range("H86").Value = 56221524
$ws.Range("I86").Value = 86541860
$ws.Range("J86").Value = 6950993
$ws.Range("K86").Value = 86541860
$ws.Range("L86").Value = 6950993
$ws.Range("M86").Value = -86540737
$ws.Range("N86").Value = -6953239

$ws.Range("H89").Value = 56221524
$ws.Range("I89").Value = 86541860
$ws.Range("J89").Value = 6950993
$ws.Range("K89").Value = 432709300
$ws.Range("L89").Value = 34754965
$ws.Range("M89").Value = -432703684
$ws.Range("N89").Value = -34766197

$ws.Range("H100").Value = 2996.875
$ws.Range("J100").Value = 4742.25
$ws.Range("L100").Value = 4742.25
$ws.Range("N100").Value = -5824.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1547739.6
$ws.Range("I32").Value = 1670556.4
$ws.Range("J32").Value = 12529.167
$ws.Range("K32").Value = 1670556.4
$ws.Range("L32").Value = 12529.167
$ws.Range("M32").Value = -1670269.4
$ws.Range("N32").Value = -13103.167

$ws.Range("H38").Value = 69000
$ws.Range("I38").Value = 69000
$ws.Range("K38").Value = 69000
$ws.Range("M38").Value = -68533

$ws.Range("H57").Value = 4527
$ws.Range("I57").Value = 4527
$ws.Range("K57").Value = 4527
$ws.Range("M57").Value = -4043

$ws.Range("H97").Value = 27778068
$ws.Range("I97").Value = 434
$ws.Range("K97").Value = 434
$ws.Range("M97").Value = 62

$ws.Range("H122").Value = 13358.579
$ws.Range("I122").Value = 13863.3125
$ws.Range("K122").Value = 41589.9375
$ws.Range("M122").Value = -39139.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6785.898
$ws.Range("I31").Value = 2392.3333
$ws.Range("K31").Value = 2392.3333
$ws.Range("M31").Value = -2097.3333

$ws.Range("H33").Value = 4999.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H34").Value = 6785.898
$ws.Range("I34").Value = 2392.3333
$ws.Range("K34").Value = 2392.3333
$ws.Range("M34").Value = -2190.3333

$ws.Range("H58").Value = 11369209
$ws.Range("I58").Value = 20001488
$ws.Range("J58").Value = 10946.263
$ws.Range("K58").Value = 20001488
$ws.Range("L58").Value = 10946.263
$ws.Range("M58").Value = -20001285
$ws.Range("N58").Value = -11352.263

$ws.Range("H62").Value = 15631258
$ws.Range("J62").Value = 6026
$ws.Range("L62").Value = 6026
$ws.Range("N62").Value = -7274

$ws.Range("H65").Value = 15631258
$ws.Range("J65").Value = 6026
$ws.Range("L65").Value = 30130
$ws.Range("N65").Value = -36370

$ws.Range("H134").Value = 6539.0293
$ws.Range("I134").Value = 1700.6666
$ws.Range("J134").Value = 9178.137000000001
$ws.Range("K134").Value = 5101.9998
$ws.Range("L134").Value = 27534.411
$ws.Range("M134").Value = -2566.9998
$ws.Range("N134").Value = -32604.411

$ws.Range("H136").Value = 11369209
$ws.Range("I136").Value = 20001488
$ws.Range("J136").Value = 10946.263
$ws.Range("K136").Value = 60004464
$ws.Range("L136").Value = 32838.789
$ws.Range("M136").Value = -60001914
$ws.Range("N136").Value = -37938.789

$ws.Range("H141").Value = 62268.85
$ws.Range("J141").Value = 63435.844
$ws.Range("L141").Value = 63435.844
$ws.Range("N141").Value = -73795.844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 366296300
$ws.Range("I4").Value = 49444444
$ws.Range("K4").Value = 148333332
$ws.Range("M4").Value = -148333220

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H37").Value = 97833
$ws.Range("J37").Value = 97833
$ws.Range("L37").Value = 293499
$ws.Range("N37").Value = -293723

$ws.Range("H131").Value = 1082.3667
$ws.Range("J131").Value = 1573
$ws.Range("L131").Value = 4719
$ws.Range("N131").Value = -14799

$ws.Range("H138").Value = 76350.78999999999
$ws.Range("J138").Value = 9012.5
$ws.Range("L138").Value = 27037.5
$ws.Range("N138").Value = -37317.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 306666660
$ws.Range("J11").Value = 450000000
$ws.Range("L11").Value = 450000000
$ws.Range("N11").Value = -450000278

$ws.Range("H58").Value = 73000
$ws.Range("J58").Value = 73000
$ws.Range("L58").Value = 73000
$ws.Range("N58").Value = -73554

$ws.Range("H70").Value = 9115.9
$ws.Range("I70").Value = 7411.4443
$ws.Range("J70").Value = 10510.454
$ws.Range("K70").Value = 7411.4443
$ws.Range("L70").Value = 10510.454
$ws.Range("M70").Value = -7141.4443
$ws.Range("N70").Value = -11050.454

$ws.Range("H73").Value = 9115.9
$ws.Range("I73").Value = 7411.4443
$ws.Range("J73").Value = 10510.454
$ws.Range("K73").Value = 7411.4443
$ws.Range("L73").Value = 10510.454
$ws.Range("M73").Value = -6475.4443
$ws.Range("N73").Value = -12382.454

$ws.Range("H97").Value = 875
$ws.Range("I97").Value = 874.5
$ws.Range("J97").Value = 876.875
$ws.Range("K97").Value = 874.5
$ws.Range("L97").Value = 876.875
$ws.Range("M97").Value = -378.5
$ws.Range("N97").Value = -1868.875

$ws.Range("H126").Value = 2669.4546
$ws.Range("J126").Value = 3277.6667
$ws.Range("L126").Value = 9833.000100000001
$ws.Range("N126").Value = -14773.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6097.2
$ws.Range("I7").Value = 4995.7617
$ws.Range("K7").Value = 4995.7617
$ws.Range("M7").Value = -4883.7617

$ws.Range("H13").Value = 871
$ws.Range("I13").Value = 499
$ws.Range("K13").Value = 499
$ws.Range("M13").Value = -359

$ws.Range("H18").Value = 55503.5
$ws.Range("I18").Value = 55503.5
$ws.Range("K18").Value = 55503.5
$ws.Range("M18").Value = -55331.5

$ws.Range("H40").Value = 6970.4414
$ws.Range("I40").Value = 5999.1665
$ws.Range("K40").Value = 5999.1665
$ws.Range("M40").Value = -5863.1665

$ws.Range("H55").Value = 34483096
$ws.Range("I55").Value = 71428640
$ws.Range("J55").Value = 583.06665
$ws.Range("K55").Value = 71428640
$ws.Range("L55").Value = 583.06665
$ws.Range("M55").Value = -71428467
$ws.Range("N55").Value = -929.06665

$ws.Range("H93").Value = 6959.9165
$ws.Range("I93").Value = 4588
$ws.Range("K93").Value = 4588
$ws.Range("M93").Value = -3340

$ws.Range("H107").Value = 2399.4
$ws.Range("I107").Value = 2399.4
$ws.Range("K107").Value = 2399.4
$ws.Range("M107").Value = -479.4000000000001

$ws.Range("H122").Value = 3822.5557
$ws.Range("I122").Value = 2839.5854
$ws.Range("J122").Value = 6922.6924
$ws.Range("K122").Value = 8518.7562
$ws.Range("L122").Value = 20768.0772
$ws.Range("M122").Value = -6068.7562
$ws.Range("N122").Value = -25668.0772

$ws.Range("H126").Value = 6097.2
$ws.Range("I126").Value = 4995.7617
$ws.Range("K126").Value = 14987.2851
$ws.Range("M126").Value = -12517.2851

$ws.Range("H132").Value = 8480427
$ws.Range("I132").Value = 17244336
$ws.Range("J132").Value = 8648.6
$ws.Range("K132").Value = 51733008
$ws.Range("L132").Value = 25945.8
$ws.Range("M132").Value = -51730478
$ws.Range("N132").Value = -31005.8

$ws.Range("H136").Value = 10058.186
$ws.Range("I136").Value = 2286.2222
$ws.Range("K136").Value = 6858.6666
$ws.Range("M136").Value = -4308.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2995
$ws.Range("I96").Value = 2995
$ws.Range("K96").Value = 2995
$ws.Range("M96").Value = -1622

$ws.Range("H100").Value = 1259.0769
$ws.Range("I100").Value = 1064.6
$ws.Range("K100").Value = 2129.2
$ws.Range("M100").Value = -1588.2
